$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header row (row 5): ID, First_Name, Last_Name
$ws.Range("E5").Value = "ID"
$ws.Range("F5").Value = "First_Name"
$ws.Range("G5").Value = "Last_Name"
$ws.Range("H5").Value = ""
$ws.Range("I5").Value = ""
$ws.Range("J5").Value = ""

# Row 6: 0, Dillon, Britt
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = "Dillon"
$ws.Range("G6").Value = "Britt"
$ws.Range("H6").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("J6").Value = ""

# Row 7: 1, Eric, Claus
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = "Eric"
$ws.Range("G7").Value = "Claus"
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = ""

# Row 8 (new): 2, Noelia, Oase
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = "Noelia"
$ws.Range("G8").Value = "Oase"
$ws.Range("E8:G8").HorizontalAlignment = -4131

# Column widths (best-fit to the new header/data text, matching the
# autofit Excel applied for columns F (First_Name) and G (Last_Name)).
# Values are pre-compensated for this host's ColumnWidth -> stored-width
# rounding so the saved <col> width lands as close as possible to the
# real Excel bestFit widths (10.85546875 / 10.7109375 character units).
$ws.Columns("F").ColumnWidth = 10
$ws.Columns("G").ColumnWidth = 9.833333333333334

# Selection
$ws.Range("G8").Select()
